$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = '6:16 PM, 08/11/2023'
$ws.Range("D2").Value = 'toothbot'
$ws.Range("E2").Value = 'Placar Libertadores do Tooth'
$ws.Range("C3").Value = '6:21 PM, 08/11/2023'
$ws.Range("E3").Value = '/placar'
$ws.Range("C4").Value = '6:21 PM, 08/11/2023'
$ws.Range("D4").Value = 'toothbot'
$ws.Range("E4").Value = 'Placar Libertadores do Tooth'
$ws.Range("C5").Value = '6:21 PM, 08/11/2023'
$ws.Range("E5").Value = 'Top Toother  |  No.Tooths'
$ws.Range("C6").Value = '6:21 PM, 08/11/2023'
$ws.Range("D6").Value = 'toothbot'
$ws.Range("E6").Value = '1.   Igor    3'
$ws.Range("C7").Value = '6:23 PM, 08/11/2023'
$ws.Range("C8").Value = '9:26 PM, 08/11/2023'
$ws.Range("D8").Value = 'toothbot'
$ws.Range("C9").Value = '9:28 PM, 08/11/2023'
$ws.Range("D9").Value = 'Igor'
$ws.Range("E9").Value = '/placar'
$ws.Range("C10").Value = '9:28 PM, 08/11/2023'
$ws.Range("E10").Value = '💩'
$ws.Range("C11").Value = '9:28 PM, 08/11/2023'
$ws.Range("D11").Value = 'Igor'
$ws.Range("E11").Value = '/placar'
$ws.Range("C12").Value = '8:22 AM, 09/11/2023'
$ws.Range("D12").Value = 'toothbot'
$ws.Range("E12").Value = 'Placar Libertadores do Tooth'
$ws.Range("C13").Value = '8:22 AM, 09/11/2023'
$ws.Range("E13").Value = 'Top Toother  |  No.Tooths'
$ws.Range("C14").Value = '8:22 AM, 09/11/2023'
$ws.Range("E14").Value = '1.   Igor           1'
$ws.Range("C15").Value = '8:23 AM, 09/11/2023'
$ws.Range("D15").Value = 'Igor'
$ws.Range("E15").Value = '💩'
$ws.Range("C16").Value = '8:24 AM, 09/11/2023'
$ws.Range("D16").Value = 'Igor'
$ws.Range("E16").Value = '/placar'
$ws.Range("C17").Value = '8:24 AM, 09/11/2023'
$ws.Range("E17").Value = 'Placar Libertadores do Tooth'
$ws.Range("C18").Value = '8:24 AM, 09/11/2023'
$ws.Range("E18").Value = 'Top Toother  |  No.Tooths'
$ws.Range("C19").Value = '8:24 AM, 09/11/2023'
$ws.Range("E19").Value = '1.   Igor           2'
$ws.Range("C20").Value = '8:27 AM, 09/11/2023'
$ws.Range("D20").Value = 'Igor'
$ws.Range("C21").Value = '8:28 AM, 09/11/2023'
$ws.Range("C22").Value = '8:28 AM, 09/11/2023'
$ws.Range("E22").Value = 'Top Toother  |  No.Tooths'
$ws.Range("C23").Value = '8:28 AM, 09/11/2023'
$ws.Range("E23").Value = '1.   Igor           5'
$ws.Range("B24").Value = 22
$ws.Range("C24").Value = '5:27 PM, 09/11/2023'
$ws.Range("B25").Value = 23
$ws.Range("C25").Value = '10:44 PM, 09/11/2023'
$ws.Range("D25").Value = 'Igor'
$ws.Range("E25").Value = '/placar'
$ws.Range("B26").Value = 24
$ws.Range("C26").Value = '10:56 PM, 09/11/2023'
$ws.Range("D26").Value = 'toothbot'
$ws.Range("E26").Value = '['
$ws.Range("B27").Value = 25
$ws.Range("C27").Value = '10:56 PM, 09/11/2023'
$ws.Range("E27").Value = '{'
$ws.Range("B28").Value = 26
$ws.Range("C28").Value = '10:56 PM, 09/11/2023'
$ws.Range("E28").Formula = '="''"'
$ws.Range("B29").Value = 27
$ws.Range("C29").Value = '10:56 PM, 09/11/2023'
$ws.Range("E29").Value = '0'
$ws.Range("B30").Value = 28
$ws.Range("C30").Value = '10:56 PM, 09/11/2023'
$ws.Range("D30").Value = 'toothbot'
$ws.Range("E30").Value = '9'
$ws.Range("B31").Value = 29
$ws.Range("C31").Value = '10:56 PM, 09/11/2023'
$ws.Range("E31").Value = '1'
$ws.Range("B32").Value = 30
$ws.Range("C32").Value = '10:56 PM, 09/11/2023'
$ws.Range("D32").Value = 'toothbot'
$ws.Range("E32").Value = '1'
$ws.Range("B33").Value = 31
$ws.Range("C33").Value = '10:56 PM, 09/11/2023'
$ws.Range("D33").Value = 'toothbot'
$ws.Range("E33").Formula = '="''"'
$ws.Range("B34").Value = 32
$ws.Range("C34").Value = '10:56 PM, 09/11/2023'
$ws.Range("D34").Value = 'toothbot'
$ws.Range("E34").Value = ':'
$ws.Range("B35").Value = 33
$ws.Range("C35").Value = '10:56 PM, 09/11/2023'
$ws.Range("E35").Value = '['
$ws.Range("B36").Value = 34
$ws.Range("C36").Value = '10:56 PM, 09/11/2023'
$ws.Range("E36").Value = '{'
$ws.Range("B37").Value = 35
$ws.Range("C37").Value = '10:56 PM, 09/11/2023'
$ws.Range("E37").Formula = '="''"'
$ws.Range("B38").Value = 36
$ws.Range("C38").Value = '10:56 PM, 09/11/2023'
$ws.Range("D38").Value = 'toothbot'
$ws.Range("E38").Value = 'I'
$ws.Range("B39").Value = 37
$ws.Range("C39").Value = '10:56 PM, 09/11/2023'
$ws.Range("D39").Value = 'toothbot'
$ws.Range("E39").Value = 'g'
$ws.Range("B40").Value = 38
$ws.Range("C40").Value = '10:56 PM, 09/11/2023'
$ws.Range("E40").Value = 'o'
$ws.Range("B41").Value = 39
$ws.Range("C41").Value = '10:56 PM, 09/11/2023'
$ws.Range("E41").Value = 'r'
$ws.Range("B42").Value = 40
$ws.Range("C42").Value = '10:56 PM, 09/11/2023'
$ws.Range("E42").Formula = '="''"'
$ws.Range("B43").Value = 41
$ws.Range("C43").Value = '10:56 PM, 09/11/2023'
$ws.Range("D43").Value = 'toothbot'
$ws.Range("E43").Value = ':'
$ws.Range("B44").Value = 42
$ws.Range("C44").Value = '10:56 PM, 09/11/2023'
$ws.Range("E44").Value = '1'
$ws.Range("B45").Value = 43
$ws.Range("C45").Value = '10:57 PM, 09/11/2023'
$ws.Range("E45").Value = '}'
$ws.Range("B46").Value = 44
$ws.Range("C46").Value = '10:57 PM, 09/11/2023'
$ws.Range("E46").Value = ']'
$ws.Range("C47").Value = '10:57 PM, 09/11/2023'
$ws.Range("D47").Value = 'toothbot'
$ws.Range("E47").Value = '}'
$ws.Range("C48").Value = '10:57 PM, 09/11/2023'
$ws.Range("D48").Value = 'toothbot'
$ws.Range("E48").Value = ']'
